$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds price values that look numeric but must stay as text
# (matches the source data which used inline strings, e.g. "38.097.40").
# Pre-format column D as Text so assigning these strings doesn't get
# silently reinterpreted as numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '38.097.40'
$ws.Range("E2").Value = '  +2.60%  '
$ws.Range("D3").Value = '2.105.63'
$ws.Range("E3").Value = '  +2.71%  '
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D5").Value = '234.71'
$ws.Range("E5").Value = '  +1.26%  '
$ws.Range("D6").Value = '0.626'
$ws.Range("E6").Value = '  +1.18%  '
$ws.Range("D7").Value = '58.17'
$ws.Range("E7").Value = '  +2.28%  '
$ws.Range("E8").Value = '  +0.02%  '
$ws.Range("D9").Value = '0.392'
$ws.Range("E9").Value = '  +2.76%  '
$ws.Range("D10").Value = '0.0780'
$ws.Range("E10").Value = '  +3.30%  '
$ws.Range("D11").Value = '0.105'
$ws.Range("E11").Value = '  +2.85%  '
$ws.Range("D12").Value = '2.406.69'
$ws.Range("E12").Value = '  +2.19%  '
$ws.Range("D13").Value = '14.63'
$ws.Range("E13").Value = '  +2.80%  '
$ws.Range("D14").Value = '21.64'
$ws.Range("E14").Value = '  +4.12%  '
$ws.Range("D15").Value = '0.784'
$ws.Range("E15").Value = '  +1.54%  '
$ws.Range("D16").Value = '5.27'
$ws.Range("E16").Value = '  +2.66%  '
$ws.Range("D17").Value = '2.109.56'
$ws.Range("E17").Value = '  +2.76%  '
$ws.Range("D18").Value = '38.009.40'
$ws.Range("E18").Value = '  +2.54%  '
$ws.Range("E19").Value = '  -2.23%  '
$ws.Range("D20").Value = '70.87'
$ws.Range("E20").Value = '  +2.54%  '
$ws.Range("D21").Value = '0.0₃0827'
$ws.Range("E21").Value = '  +2.64%  '
$ws.Range("D22").Value = '228.43'
$ws.Range("E22").Value = '  +1.61%  '
$ws.Range("D23").Value = '1.00'
$ws.Range("E23").Value = '  +0.01%  '
$ws.Range("D24").Value = '2.41'
$ws.Range("E24").Value = '  +0.75%  '
$ws.Range("E25").Value = '  +1.88%  '
$ws.Range("D26").Value = '168.19'
$ws.Range("E26").Value = '  +1.40%  '
$ws.Range("D27").Value = '0.139'
$ws.Range("E27").Value = '  +10.85%  '
$ws.Range("D28").Value = '9.03'
$ws.Range("E28").Value = '  +3.37%  '
$ws.Range("E29").Value = '  -0.61%  '
$ws.Range("D30").Value = '19.51'
$ws.Range("E30").Value = '  +3.01%  '
$ws.Range("E31").Value = '  +1.64%  '
$ws.Range("D32").Value = '4.68'
$ws.Range("E32").Value = '  +5.61%  '
$ws.Range("E33").Value = '  +4.57%  '
$ws.Range("D34").Value = '0.0627'
$ws.Range("E34").Value = '  +2.14%  '
$ws.Range("D35").Value = '4.63'
$ws.Range("E35").Value = '  +1.46%  '
$ws.Range("D36").Value = '3.47'
$ws.Range("E36").Value = '  +6.64%  '
$ws.Range("E37").Value = '  +5.00%  '
$ws.Range("E38").Value = '  +0.08%  '
$ws.Range("D39").Value = '5.48'
$ws.Range("E39").Value = '  -3.92%  '
$ws.Range("D40").Value = '0.0992'
$ws.Range("E40").Value = '  +6.79%  '
$ws.Range("D42").Value = '97.18'
$ws.Range("E42").Value = '  +0.89%  '
$ws.Range("E43").Value = '  +3.32%  '
$ws.Range("D44").Value = '1.458.56'
$ws.Range("E44").Value = '  -1.33%  '
$ws.Range("E45").Value = '  -0.77%  '
$ws.Range("E46").Value = '  +7.39%  '
$ws.Range("D47").Value = '4.15'
$ws.Range("E47").Value = '  -5.18%  '
$ws.Range("D48").Value = '1.06'
$ws.Range("E48").Value = '  +4.86%  '
$ws.Range("B49").Value = 'FraxShare'
$ws.Range("C49").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D49").Value = '7.32'
$ws.Range("E49").Value = '  +2.17%  '
$ws.Range("B50").Value = 'MXToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D50").Value = '3.03'
$ws.Range("E50").Value = '  +3.28%  '
$ws.Range("D51").Value = '2.301.82'
$ws.Range("E51").Value = '  +2.65%  '
